$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("B1").Value = "kostnadd (kr)"

# Row 2 stays the same: domene / 150 (unchanged)
$ws.Range("A2").Value = "domene"
$ws.Range("B2").Value = 150

# New row order / values
$ws.Range("A3").Value = "lisenser"
$ws.Range("B3").Value = 1500

$ws.Range("A4").Value = "arbeid for prototype"
$ws.Range("B4").Value = 2000

$ws.Range("A5").Value = "arbeid for og ferdigstille"
$ws.Range("B5").Value = 2000

$ws.Range("A6").Value = "SUM"
$ws.Range("B6").Formula = "=SUM(B2:B5)"

$ws.Range("F9").Select()
